$wb = $excel.ActiveWorkbook

# --- Sheet "About": add Notes section below existing content ---
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("A9").Value = "Notes:"
$wsAbout.Range("A10").Value = "This policy covers improvements in air sealing, framing, and insulation, so it applies to the "
$wsAbout.Range("A11").Value = """envelope"" component."
$wsAbout.Range("A11").Select()

# --- Sheet "PPEIdtICEaT": update header row text, wrap, and row height ---
$wsMain = $wb.Worksheets.Item("PPEIdtICEaT")
$wsMain.Range("A1").Value = "Building Component Efficiency Improvement (dimensionless)"
$wsMain.Range("A1").WrapText = $true
$wsMain.Range("B1").WrapText = $true
$wsMain.Rows.Item(1).RowHeight = 45
